# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect newly output data, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1089
$wsExhibit.Range("F5").Value = 3585
$wsExhibit.Range("F10").Value = 20
$wsExhibit.Range("F11").Value = 130
$wsExhibit.Range("F13").Value = 247
$wsExhibit.Range("F14").Value = 48
$wsExhibit.Range("F15").Value = 94
$wsExhibit.Range("F16").Value = 2813
$wsExhibit.Range("F17").Value = 1132

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1089
$wsAll.Range("F6").Value = 3586
$wsAll.Range("F12").Value = 20
$wsAll.Range("F13").Value = 130
$wsAll.Range("F15").Value = 249
$wsAll.Range("F16").Value = 48
$wsAll.Range("F17").Value = 94
$wsAll.Range("F18").Value = 2813
$wsAll.Range("F19").Value = 1132

$wb.Save()
